$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 999
$ws.Range("I4").Value = 999
$ws.Range("K4").Value = 999
$ws.Range("M4").Value = -885
$ws.Range("H6").Value = 200018
$ws.Range("I6").Value = 200018
$ws.Range("K6").Value = 600054
$ws.Range("M6").Value = -599942
$ws.Range("H19").Value = 1170.2222
$ws.Range("I19").Value = 1327.5
$ws.Range("J19").Value = 855.6667
$ws.Range("K19").Value = 1327.5
$ws.Range("L19").Value = 855.6667
$ws.Range("M19").Value = -1152.5
$ws.Range("N19").Value = -1205.6667
$ws.Range("H26").Value = 2900
$ws.Range("I26").Value = 2900
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 2900
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -2556
$ws.Range("N26").ClearContents()
$ws.Range("H33").Value = 323.2857
$ws.Range("J33").Value = 400
$ws.Range("L33").Value = 400
$ws.Range("N33").Value = -858
$ws.Range("H62").Value = 1999.5
$ws.Range("I62").Value = 1999.5
$ws.Range("K62").Value = 1999.5
$ws.Range("M62").Value = -1375.5
$ws.Range("H65").Value = 1999.5
$ws.Range("I65").Value = 1999.5
$ws.Range("K65").Value = 9997.5
$ws.Range("M65").Value = -6877.5
$ws.Range("H105").Value = 35335
$ws.Range("J105").Value = 35335
$ws.Range("L105").Value = 35335
$ws.Range("N105").Value = -42323
$ws.Range("H111").Value = 3724.6667
$ws.Range("I111").Value = 3724.6667
$ws.Range("K111").Value = 11174.0001
$ws.Range("M111").Value = -8107.000100000001
$ws.Range("H141").Value = 5240.7144
$ws.Range("I141").Value = 5831.8335
$ws.Range("J141").Value = 1694
$ws.Range("K141").Value = 17495.5005
$ws.Range("L141").Value = 5082
$ws.Range("M141").Value = -12315.5005
$ws.Range("N141").Value = -15442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9076.166999999999
$ws.Range("I32").Value = 8433.588
$ws.Range("K32").Value = 8433.588
$ws.Range("M32").Value = -8146.588
$ws.Range("H45").Value = 2533
$ws.Range("I45").Value = 2388.3333
$ws.Range("K45").Value = 2388.3333
$ws.Range("M45").Value = -2011.3333
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H110").Value = 6966.0713
$ws.Range("I110").Value = 7776.7646
$ws.Range("K110").Value = 7776.7646
$ws.Range("M110").Value = -5731.7646
$ws.Range("H132").Value = 3608.6128
$ws.Range("I132").Value = 3815.9048
$ws.Range("K132").Value = 11447.7144
$ws.Range("M132").Value = -8917.714399999999
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 65421.25
$ws.Range("I88").Value = 26000
$ws.Range("J88").Value = 78561.664
$ws.Range("K88").Value = 26000
$ws.Range("L88").Value = 78561.664
$ws.Range("M88").Value = -25594
$ws.Range("N88").Value = -79373.664
$ws.Range("H91").Value = 65421.25
$ws.Range("I91").Value = 26000
$ws.Range("J91").Value = 78561.664
$ws.Range("K91").Value = 26000
$ws.Range("L91").Value = 78561.664
$ws.Range("M91").Value = -24596
$ws.Range("N91").Value = -81369.664
$ws.Range("H105").Value = 2376.75
$ws.Range("I105").Value = 2169
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2169
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -422
$ws.Range("N105").Value = -6494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4446534
$ws.Range("I22").Value = 2200.5
$ws.Range("J22").Value = 8002000.5
$ws.Range("K22").Value = 2200.5
$ws.Range("L22").Value = 8002000.5
$ws.Range("M22").Value = -1850.5
$ws.Range("N22").Value = -8002700.5
$ws.Range("H31").Value = 3142.111
$ws.Range("I31").Value = 2570.8
$ws.Range("J31").Value = 3856.25
$ws.Range("K31").Value = 2570.8
$ws.Range("L31").Value = 3856.25
$ws.Range("M31").Value = -2275.8
$ws.Range("N31").Value = -4446.25
$ws.Range("H34").Value = 3142.111
$ws.Range("I34").Value = 2570.8
$ws.Range("J34").Value = 3856.25
$ws.Range("K34").Value = 2570.8
$ws.Range("L34").Value = 3856.25
$ws.Range("M34").Value = -2368.8
$ws.Range("N34").Value = -4260.25
$ws.Range("H58").Value = 3627.75
$ws.Range("I58").Value = 1503.6666
$ws.Range("K58").Value = 1503.6666
$ws.Range("M58").Value = -1300.6666
$ws.Range("H88").Value = 50000
$ws.Range("J88").Value = 50000
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50812
$ws.Range("H91").Value = 50000
$ws.Range("J91").Value = 50000
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -52808
$ws.Range("H132").Value = 2893.348
$ws.Range("I132").Value = 2824.8125
$ws.Range("K132").Value = 8474.4375
$ws.Range("M132").Value = -5944.4375
$ws.Range("H134").Value = 2512.6
$ws.Range("J134").Value = 2649.4
$ws.Range("L134").Value = 7948.200000000001
$ws.Range("N134").Value = -13018.2
$ws.Range("H136").Value = 3627.75
$ws.Range("I136").Value = 1503.6666
$ws.Range("K136").Value = 4510.9998
$ws.Range("M136").Value = -1960.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3399.5
$ws.Range("I3").Value = 3399.5
$ws.Range("K3").Value = 10198.5
$ws.Range("M3").Value = -10086.5
$ws.Range("H7").Value = 86746.164
$ws.Range("I7").Value = 142880.72
$ws.Range("K7").Value = 428642.16
$ws.Range("M7").Value = -428530.16
$ws.Range("H17").Value = 2131.6667
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 2710.7144
$ws.Range("K17").Value = 315
$ws.Range("L17").Value = 8132.1432
$ws.Range("M17").Value = -146
$ws.Range("N17").Value = -8470.143199999999
$ws.Range("H49").Value = 2708
$ws.Range("J49").Value = 2437.5
$ws.Range("L49").Value = 7312.5
$ws.Range("N49").Value = -7624.5
$ws.Range("H75").Value = 784
$ws.Range("J75").Value = 782
$ws.Range("L75").Value = 2346
$ws.Range("N75").Value = -4342
$ws.Range("H78").Value = 784
$ws.Range("J78").Value = 782
$ws.Range("L78").Value = 7038
$ws.Range("N78").Value = -17022
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H81").Value = 1474.75
$ws.Range("I81").Value = 1500
$ws.Range("K81").Value = 4500
$ws.Range("M81").Value = -3377
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H84").Value = 1474.75
$ws.Range("I84").Value = 1500
$ws.Range("K84").Value = 13500
$ws.Range("M84").Value = -7884
$ws.Range("H139").Value = 4290.5
$ws.Range("I139").Value = 776.3333
$ws.Range("K139").Value = 2328.9999
$ws.Range("M139").Value = 2811.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 8600.4
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10232
$ws.Range("H97").Value = 2171.25
$ws.Range("I97").Value = 2145
$ws.Range("K97").Value = 2145
$ws.Range("M97").Value = -1649
$ws.Range("H132").Value = 3028.4348
$ws.Range("I132").Value = 2892.4736
$ws.Range("K132").Value = 8677.4208
$ws.Range("M132").Value = -6147.4208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H7").Value = 1966.1111
$ws.Range("I7").Value = 1966.1111
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1966.1111
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1854.1111
$ws.Range("N7").ClearContents()
$ws.Range("H82").Value = 1355.0588
$ws.Range("I82").Value = 554.7778
$ws.Range("J82").Value = 2255.375
$ws.Range("K82").Value = 554.7778
$ws.Range("L82").Value = 2255.375
$ws.Range("M82").Value = -193.7778
$ws.Range("N82").Value = -2977.375
$ws.Range("H85").Value = 1355.0588
$ws.Range("I85").Value = 554.7778
$ws.Range("J85").Value = 2255.375
$ws.Range("K85").Value = 554.7778
$ws.Range("L85").Value = 2255.375
$ws.Range("M85").Value = 693.2222
$ws.Range("N85").Value = -4751.375
$ws.Range("H100").Value = 1960.6666
$ws.Range("I100").Value = 1963.8572
$ws.Range("K100").Value = 1963.8572
$ws.Range("M100").Value = -1422.8572
$ws.Range("H126").Value = 1966.1111
$ws.Range("I126").Value = 1966.1111
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5898.3333
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3428.3333
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 34000000
$ws.Range("I2").Value = 100000000
$ws.Range("K2").Value = 100000000
$ws.Range("M2").Value = -99999888
$ws.Range("H4").Value = 5736000
$ws.Range("J4").Value = 5736000
$ws.Range("L4").Value = 5736000
$ws.Range("N4").Value = -5736226
$ws.Range("H132").Value = 3560.2
$ws.Range("I132").Value = 2791.5
$ws.Range("K132").Value = 8374.5
$ws.Range("M132").Value = -5844.5
